$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A46").Value = "Golang Engineer"
$ws.Range("B46").Value = "https://www.dice.com/job-detail/dedc818b-64d3-44f3-98ae-c16edc7d047d"
$ws.Range("C46").Value = "Phoenix, Arizona"
$ws.Range("D46").Value = "Contract"
$ws.Range("E46").Value = "USD 140,000.00 - 150,000.00 per year"
$ws.Range("F46").Value = "HMG America"
